# invitation task: correct response added in schedule for module 2b
#
# Adds a new column I ("correct") holding the letter of the correct
# forced-choice response (a/b/c/d) for every trial row, and fixes the
# ordering of a handful of E:H answer-option cells that had been entered
# out of order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: header + per-row correct answer -----------------------
$correct = @{
    1  = "correct"
    2  = "d"
    3  = "c"
    4  = "c"
    5  = "a"
    6  = "a"
    7  = "c"
    8  = "b"
    9  = "a"
    10 = "b"
    11 = "a"
    12 = "c"
    13 = "b"
    14 = "d"
    15 = "b"
    16 = "c"
    17 = "a"
    18 = "b"
    19 = "c"
    20 = "a"
    21 = "c"
    22 = "a"
    23 = "a"
    24 = "d"
    25 = "a"
    26 = "c"
    27 = "d"
    28 = "b"
    29 = "a"
}

for ($row = 1; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 9)  # column I
    $cell.Value = $correct[$row]
    $cell.NumberFormat = "@"
}

# --- Fix a few rows whose answer-option cells were out of order ----------

# Row 4: F/G/H were entered rotated by one; restore the intended order.
$ws.Range("F4").Value = "5-10"
$ws.Range("G4").Value = "0-20"
$ws.Range("H4").Value = "0-5"

# Row 12: E and G were swapped.
$ws.Range("E12").Value = "5-5"
$ws.Range("G12").Value = "10-0"

# Row 15: E and F were swapped.
$ws.Range("E15").Value = "10-0"
$ws.Range("F15").Value = "0-5"

# Row 27: G and H were swapped.
$ws.Range("G27").Value = "5-10"
$ws.Range("H27").Value = "10-10"

# --- Cosmetic: keep the last-active-cell selection in sync ---------------
$ws.Range("AF25").Select()
